$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text
$ws.Range("B1").Value = "Sorted Member List as of 02/11/2021"

# Overwrite the first two data rows (currently KANTA @ row3, Narumi @ row4)
# with the two brand-new members, then insert one more new row, keep KANTA,
# and drop the remaining old members (SHIA, YUKARI, Mariko).

# Row 3 -> TOKO
$ws.Range("B3").Value = 19714
$ws.Range("C3").Value = "TOKO"
$ws.Range("D3").Value = "morimura "
$ws.Range("E3").Value = "tokoko0415@gmail.com"
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = "02-23-2021  03:18:40 PM"

# Row 4 -> KIZUKU IW
$ws.Range("B4").Value = 19283
$ws.Range("C4").Value = "KIZUKU IW"
$ws.Range("D4").Value = "yamashita"
$ws.Range("E4").Value = "ykizuku+1@gmail.com"
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = "02-19-2021  11:09:36 PM"

# Row 5 -> SHIGEKI
$ws.Range("B5").Value = 16400
$ws.Range("C5").Value = "SHIGEKI"
$ws.Range("D5").Value = "kudo"
$ws.Range("E5").Value = "sigeki1968@gmail.com"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "02-19-2021  08:55:18 PM"

# Row 6 -> KANTA (unchanged data, but credits 1 -> 1, keep id 19509)
$ws.Range("B6").Value = 19509
$ws.Range("C6").Value = "KANTA"
$ws.Range("D6").Value = "imori"
$ws.Range("E6").Value = "w-deco@joetsu.ne.jp"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "02-09-2021  11:29:32 PM"

# Remove the now-obsolete trailing row (was row 7, Mariko)
$ws.Rows.Item(7).Delete()

[void]$ws.Range("B6:G6").Select()
